$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SOL row (row 4): Anzahl 33.54 -> 42.31, Kaufpreis 21 -> 37.38
$ws.Range("B4").Value = 42.31
$ws.Range("C4").Value = 37.38

# Update ATOM row (row 6): Anzahl 1481 -> 1488
$ws.Range("B6").Value = 1488

# Remove the ADA row (row 10) entirely, shifting subsequent rows up
$ws.Rows(10).Delete()

# Update the active selection to match the saved view state
$ws.Range("B6").Select()
